$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Reorder the "Periodo Mora" values in column E (rows 16-20) from
# descending (1708,1707,1706,1705,1704) to ascending (1704..1708).
$ws.Range("E16").Value = "1704"
$ws.Range("E17").Value = "1705"
$ws.Range("E18").Value = "1706"
$ws.Range("E19").Value = "1707"
$ws.Range("E20").Value = "1708"

# Update "Salario Basico" values in column G (rows 16-20) to the new amount.
$ws.Range("G16").Value = 781242
$ws.Range("G17").Value = 781242
$ws.Range("G18").Value = 781242
$ws.Range("G19").Value = 781242
$ws.Range("G20").Value = 781242
